$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) miscPages: "about"/"accessibility"/"careers"/"privacy" URLs are moved
#    from staging.lexus.com to the CEPO proxy host, and the two
#    privacy/legal-terms + privacy/online-statement rows are removed.
# ---------------------------------------------------------------------------
$misc = $wb.Worksheets.Item("miscPages")

$misc.Range("A2").Value = "https://stg-lcom.cdn.cepo-proxy.tms.aws.lexus.com/about"
$misc.Range("A3").Value = "https://stg-lcom.cdn.cepo-proxy.tms.aws.lexus.com/about/manufacturing"
$misc.Range("A4").Value = "https://stg-lcom.cdn.cepo-proxy.tms.aws.lexus.com/about/technology"
$misc.Range("A5").Value = "https://stg-lcom.cdn.cepo-proxy.tms.aws.lexus.com/about/environment"
$misc.Range("A6").Value = "https://stg-lcom.cdn.cepo-proxy.tms.aws.lexus.com/about/philanthropy"
$misc.Range("A8").Value = "https://stg-lcom.cdn.cepo-proxy.tms.aws.lexus.com/accessibility"
$misc.Range("A9").Value = "https://stg-lcom.cdn.cepo-proxy.tms.aws.lexus.com/careers"
$misc.Range("A10").Value = "https://stg-lcom.cdn.cepo-proxy.tms.aws.lexus.com/privacy"

# Remove the two "privacy/legal-terms" and "privacy/online-statement" rows
# (rows 11 and 12), shifting everything below up.
$misc.Range("A11:A12").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp) | Out-Null

# ---------------------------------------------------------------------------
# 2) ComparePages loses its "tabSelected" flag (handled automatically by
#    activating miscPages at the end of this script, which deactivates the
#    previously active tab).
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 3) DealerPages: all dealer URLs move from staging.lexus.com to the CEPO
#    proxy host and become real hyperlinks with the Hyperlink style.
# ---------------------------------------------------------------------------
$dealers = $wb.Worksheets.Item("DealerPages")

$dealerUrls = @(
    "https://stg-lcom.cdn.cepo-proxy.tms.aws.lexus.com/dealers",
    "https://stg-lcom.cdn.cepo-proxy.tms.aws.lexus.com/dealers/63110-lexus-of-manhattan",
    "https://stg-lcom.cdn.cepo-proxy.tms.aws.lexus.com/dealers/64204-sewell-lexus",
    "https://stg-lcom.cdn.cepo-proxy.tms.aws.lexus.com/dealers/61230-mcgrath-lexus-of-chicago",
    "https://stg-lcom.cdn.cepo-proxy.tms.aws.lexus.com/dealers/60438-jim-falk-lexus-of-beverly-hills",
    "https://stg-lcom.cdn.cepo-proxy.tms.aws.lexus.com/dealers/60419-keyes-lexus",
    "https://stg-lcom.cdn.cepo-proxy.tms.aws.lexus.com/dealers/60406-lexus-santa-monica",
    "https://stg-lcom.cdn.cepo-proxy.tms.aws.lexus.com/dealers/60445-lexus-of-cerritos"
)

for ($i = 0; $i -lt $dealerUrls.Length; $i++) {
    $row = $i + 1
    $cell = $dealers.Range("A$row")
    $cell.Value = $dealerUrls[$i]
    $dealers.Hyperlinks.Add($cell, $dealerUrls[$i]) | Out-Null
    $cell.Style = "Hyperlink"
}

$dealers.Range("B12").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4) FCVPages: add the new "LF-30 Electrified" concept page as a fifth row
#    with its own hyperlink.
# ---------------------------------------------------------------------------
$fcv = $wb.Worksheets.Item("FCVPages")

$fcvUrl = "https://stg-lcom.cdn.cepo-proxy.tms.aws.lexus.com/concept/LF-30-Electrified"
$fcvCell = $fcv.Range("A5")
$fcvCell.Value = $fcvUrl
$fcv.Hyperlinks.Add($fcvCell, $fcvUrl) | Out-Null
$fcvCell.Style = "Hyperlink"

$fcv.Range("A5").Select() | Out-Null

# ---------------------------------------------------------------------------
# 5) miscPages becomes the active / selected sheet (activeTab="1" in the
#    workbook view). This must happen last so it "wins" as the final
#    activation recorded by the workbook.
# ---------------------------------------------------------------------------
$misc.Activate()
$misc.Range("C32").Select() | Out-Null
